$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rotate C1/D1/E1 values ---
# Before: C1="max", D1="prediction", E1="rejection-f"
# After:  C1="prediction", D1="rejection-f", E1="max"
$ws.Cells.Item(1, 3).Value = "prediction"
$ws.Cells.Item(1, 4).Value = "rejection-f"
$ws.Cells.Item(1, 5).Value = "max"

# --- Data rows (2 through 56) ---
# For each row: column C becomes the same text as column D
# ("o__Coriobacteriales"), and column E becomes the number 1.
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "o__Coriobacteriales"
    $ws.Cells.Item($r, 5).Value = 1
}
